$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("使用案例")

$ws.Range("A28").Value = 29
$ws.Range("B28").Value = "1:年度編號"
$ws.Range("C28").Value = "L6"
$ws.Range("D28").Value = "RvNo"
$ws.Range("E28").Value = 999999
$ws.Range("F28").Value = "銷帳編號"
$ws.Range("G28").Value = 'gSeqCom.getSeqNo(this.getTxBuffer().getMgBizDate().getTbsDy(), 1, "L6", "RvNo", 999999, titaVo),'

$ws.Range("M29").Value = "6)"
